$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6409899155710832
$ws.Range("C2").Value = 0.1702344879934543
$ws.Range("E2").Value = 0.1166335526239486
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.002483308034040507
$ws.Range("K2").Value = 0.3259199722913024
$ws.Range("L2").Value = 0.1941625722818543
$ws.Range("O2").Value = 3.896598061147927
$ws.Range("B3").Value = 0.5986403686690664
$ws.Range("C3").Value = 0.170822735542334
$ws.Range("E3").Value = 0.1161758985404262
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.002485766210063481
$ws.Range("K3").Value = 0.2900583896036153
$ws.Range("L3").Value = 0.1870998553846164
$ws.Range("O3").Value = 3.937184168298302
$ws.Range("B4").Value = 0.5728138058209424
$ws.Range("C4").Value = 0.1712176997421793
$ws.Range("E4").Value = 0.1159574664341783
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.002487355602100291
$ws.Range("K4").Value = 0.2680360249202209
$ws.Range("L4").Value = 0.1828632342665202
$ws.Range("O4").Value = 3.964266474408944
$ws.Range("B5").Value = 0.5623342598426007
$ws.Range("C5").Value = 0.1713871752353278
$ws.Range("E5").Value = 0.1158842047281752
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002488023483792933
$ws.Range("K5").Value = 0.2590614486809386
$ws.Range("L5").Value = 0.1811619752809435
$ws.Range("O5").Value = 3.975846238703127
$ws.Range("B6").Value = 0.5605968762783675
$ws.Range("C6").Value = 0.1714158322684796
$ws.Range("E6").Value = 0.1158729915262349
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002488135606140981
$ws.Range("K6").Value = 0.2575712254433142
$ws.Range("L6").Value = 0.180881006635687
$ws.Range("O6").Value = 3.977801869765756
$ws.Range("B7").Value = 0.5726722919486349
$ws.Range("C7").Value = 0.171219950792981
$ws.Range("E7").Value = 0.115956414601925
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002487364527359759
$ws.Range("K7").Value = 0.2679149910999286
$ws.Range("L7").Value = 0.182840188365887
$ws.Range("O7").Value = 3.96442044282152
$ws.Range("B8").Value = 0.6263516044737401
$ws.Range("C8").Value = 0.1704303234606179
$ws.Range("E8").Value = 0.1164627773236226
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.002484139032095745
$ws.Range("K8").Value = 0.3135558661634832
$ws.Range("L8").Value = 0.191706651689401
$ws.Range("O8").Value = 3.910143242449564
$ws.Range("B9").Value = 0.732990974277925
$ws.Range("C9").Value = 0.16914854717335
$ws.Range("E9").Value = 0.1179516435960579
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002478446387105329
$ws.Range("K9").Value = 0.4030135338434491
$ws.Range("L9").Value = 0.2098844040383199
$ws.Range("O9").Value = 3.82087234488074
$ws.Range("B10").Value = 0.8121530423523495
$ws.Range("C10").Value = 0.1683676344279377
$ws.Range("E10").Value = 0.1193473572485111
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002474645803309712
$ws.Range("K10").Value = 0.4686930233287399
$ws.Range("L10").Value = 0.2237203555326914
$ws.Range("O10").Value = 3.765764031698012
$ws.Range("B11").Value = 0.848338212128624
$ws.Range("C11").Value = 0.1680469268237985
$ws.Range("E11").Value = 0.1200477665634772
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002472998912748034
$ws.Range("K11").Value = 0.4985591055364864
$ws.Range("L11").Value = 0.2301189182709606
$ws.Range("O11").Value = 3.742972317160564
$ws.Range("B12").Value = 0.8620650490538821
$ws.Range("C12").Value = 0.1679304194312579
$ws.Range("E12").Value = 0.1203223982203809
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002472387010935562
$ws.Range("K12").Value = 0.5098664776598127
$ws.Range("L12").Value = 0.2325568667579319
$ws.Range("O12").Value = 3.734669447438648
$ws.Range("B13").Value = 0.8591076597920733
$ws.Range("C13").Value = 0.1679552922279086
$ws.Range("E13").Value = 0.1202628334828333
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.00247251827368366
$ws.Range("K13").Value = 0.507431342404999
$ws.Range("L13").Value = 0.2320311471968637
$ws.Range("O13").Value = 3.736443031370072
$ws.Range("B14").Value = 0.8494670443727159
$ws.Range("C14").Value = 0.1680372428958634
$ws.Range("E14").Value = 0.1200701723391511
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002472948336385006
$ws.Range("K14").Value = 0.4994894184462453
$ws.Range("L14").Value = 0.2303191905547237
$ws.Range("O14").Value = 3.742282662133505
$ws.Range("B15").Value = 0.8435650322641948
$ws.Range("C15").Value = 0.168088082234064
$ws.Range("E15").Value = 0.1199533857522539
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002473213289407081
$ws.Range("K15").Value = 0.4946244491161451
$ws.Range("L15").Value = 0.2292725125061423
$ws.Range("O15").Value = 3.745902312650713
$ws.Range("B16").Value = 0.8097916954943116
$ws.Range("C16").Value = 0.1683892859022933
$ws.Range("E16").Value = 0.1193029003703145
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002474755077243134
$ws.Range("K16").Value = 0.4667409181657831
$ws.Range("L16").Value = 0.2233042909294483
$ws.Range("O16").Value = 3.76729938211659
$ws.Range("B17").Value = 0.7891168898499927
$ws.Range("C17").Value = 0.1685828892476664
$ws.Range("E17").Value = 0.1189206115298163
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.00247572188125308
$ws.Range("K17").Value = 0.4496318403745079
$ws.Range("L17").Value = 0.2196696946258072
$ws.Range("O17").Value = 3.781009341483525
$ws.Range("B18").Value = 0.7772417038354433
$ws.Range("C18").Value = 0.1686974969004567
$ws.Range("E18").Value = 0.1187068944872571
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002476285684084049
$ws.Range("K18").Value = 0.4397900696295665
$ws.Range("L18").Value = 0.217589013463666
$ws.Range("O18").Value = 3.789109278176156
$ws.Range("B19").Value = 0.7732238110849039
$ws.Range("C19").Value = 0.1687368605031025
$ws.Range("E19").Value = 0.1186355929610379
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.00247647790599887
$ws.Range("K19").Value = 0.4364576465288508
$ws.Range("L19").Value = 0.2168862234078546
$ws.Range("O19").Value = 3.791888577488393
$ws.Range("B20").Value = 0.791316064467054
$ws.Range("C20").Value = 0.16856194346542
$ws.Range("E20").Value = 0.1189606688497555
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002475618164364656
$ws.Range("K20").Value = 0.4514532469434016
$ws.Range("L20").Value = 0.2200555855110906
$ws.Range("O20").Value = 3.779527707722849
$ws.Range("B21").Value = 0.852298073044409
$ws.Range("C21").Value = 0.1680130382153422
$ws.Range("E21").Value = 0.1201265065504131
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002472821698390287
$ws.Range("K21").Value = 0.5018222199605873
$ws.Range("L21").Value = 0.2308216285580471
$ws.Range("O21").Value = 3.740558520024933
$ws.Range("B22").Value = 0.8922946410729651
$ws.Range("C22").Value = 0.1676830642150691
$ws.Range("E22").Value = 0.1209432426425288
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002471062448535537
$ws.Range("K22").Value = 0.534727794716332
$ws.Range("L22").Value = 0.237944973361806
$ws.Range("O22").Value = 3.717001086326889
$ws.Range("B23").Value = 0.8709350265199305
$ws.Range("C23").Value = 0.1678565545945148
$ws.Range("E23").Value = 0.1205023268452656
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002471995152206666
$ws.Range("K23").Value = 0.5171668866162236
$ws.Range("L23").Value = 0.2341351651299988
$ws.Range("O23").Value = 3.729399133108245
$ws.Range("B24").Value = 0.7903217831960205
$ws.Range("C24").Value = 0.1685714027635754
$ws.Range("E24").Value = 0.1189425400399671
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002475665030126361
$ws.Range("K24").Value = 0.4506298062299265
$ws.Range("L24").Value = 0.2198810965250573
$ws.Range("O24").Value = 3.78019687562815
$ws.Range("B25").Value = 0.7039975272757601
$ws.Range("C25").Value = 0.1694669355625926
$ws.Range("E25").Value = 0.1174958250837165
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002479919075965318
$ws.Range("K25").Value = 0.3788194696261087
$ws.Range("L25").Value = 0.2048823245726936
$ws.Range("O25").Value = 3.843183263851145
